$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple per-row updates: Price (D) and/or Volume(1h) (E) changes ---
# Numeric-looking Price strings are prefixed with a leading apostrophe so
# Excel stores them verbatim as text (matching the source data's exact
# formatting) instead of auto-converting / trimming them as numbers.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "57.885.04"
$ws.Range("E2").Value = "  -0.90%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.566.48"
$ws.Range("E3").Value = "  -3.12%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.03%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'515.77"
$ws.Range("E5").Value = "  -1.25%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'142.37"
$ws.Range("E6").Value = "  -1.86%  "

# Row 7 - USDC (unchanged)

# Row 8 - XRP
$ws.Range("E8").Value = "  -1.24%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.580.99"
$ws.Range("E9").Value = "  -2.72%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  -2.89%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -2.07%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "'0.323"
$ws.Range("E12").Value = "  -4.77%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  -1.15%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.017.21"
$ws.Range("E14").Value = "  -3.17%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "57.866.39"
$ws.Range("E15").Value = "  -0.93%  "

# Row 16 - Avalanche
$ws.Range("D16").Value = "'20.24"
$ws.Range("E16").Value = "  -3.53%  "

# --- Rows 17/18 swap: ShibaInu <-> WrappedEther (with updated values) ---
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.622.46"
$ws.Range("E17").Value = "  -0.88%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.0000133"
$ws.Range("E18").Value = "  -2.45%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "'339.16"
$ws.Range("E19").Value = "  +0.01%  "

# Row 20 - Polkadot
$ws.Range("E20").Value = "  -2.61%  "

# Row 21 - Chainlink
$ws.Range("E21").Value = "  -2.54%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -0.36%  "

# Row 23 - Dai
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.36%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'65.39"
$ws.Range("E24").Value = "  +1.50%  "

# Row 25 - Kaspa
$ws.Range("E25").Value = "  -0.69%  "

# --- Rows 26/27 swap: Polygon <-> Binance-PegBSC-USD (with updated values) ---
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("B27").Value = "Polygon"
$ws.Range("C27").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D27").Value = "'0.400"
$ws.Range("E27").Value = "  -6.07%  "

# Row 28 - WrappedeETH
$ws.Range("D28").Value = "2.687.92"
$ws.Range("E28").Value = "  -2.89%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("E29").Value = "  -2.59%  "

# Row 30 - PEPE
$ws.Range("D30").Value = "0.0₃0741"
$ws.Range("E30").Value = "  -7.33%  "

# Row 31 - USDe
$ws.Range("E31").Value = "  -0.07%  "

# Row 32 - Aptos
$ws.Range("D32").Value = "'6.27"
$ws.Range("E32").Value = "  -6.40%  "

# Row 33 - PancakeSwap
$ws.Range("E33").Value = "  -1.28%  "

# Row 34 - EthereumClassic
$ws.Range("D34").Value = "'18.63"
$ws.Range("E34").Value = "  -1.33%  "

# Row 35 - Monero
$ws.Range("D35").Value = "'149.80"
$ws.Range("E35").Value = "  -1.81%  "

# Row 36 - NEARProtocol
$ws.Range("D36").Value = "'3.98"
$ws.Range("E36").Value = "  -4.20%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  -3.97%  "

# Row 38 - SuiNetwork
$ws.Range("D38").Value = "'0.872"
$ws.Range("E38").Value = "  -4.30%  "

# Row 39 - OKB
$ws.Range("D39").Value = "'36.03"
$ws.Range("E39").Value = "  -2.10%  "

# Row 40 - Stacks
$ws.Range("E40").Value = "  -0.49%  "

# Row 41 - Fetch.AI
$ws.Range("D41").Value = "'0.830"
$ws.Range("E41").Value = "  -4.68%  "

# Row 42 - Filecoin
$ws.Range("E42").Value = "  -3.07%  "

# Row 43 - FirstDigitalUSD
$ws.Range("D43").Value = "'0.998"
$ws.Range("E43").Value = "  -0.18%  "

# Row 44 - Bittensor
$ws.Range("D44").Value = "'269.62"
$ws.Range("E44").Value = "  -1.79%  "

# Row 45 - WhiteBITCoin
$ws.Range("D45").Value = "'10.67"
$ws.Range("E45").Value = "  +0.43%  "

# Row 46 - Stellar
$ws.Range("D46").Value = "'0.0951"
$ws.Range("E46").Value = "  -2.11%  "

# Row 47 - Mantle
$ws.Range("E47").Value = "  -3.94%  "

# Row 48 - EnergySwap
$ws.Range("D48").Value = "'18.68"
$ws.Range("E48").Value = "  -4.02%  "

# Row 49 - Hedera
$ws.Range("D49").Value = "'0.0520"
$ws.Range("E49").Value = "  -3.11%  "

# Row 50 - Maker
$ws.Range("D50").Value = "1.979.57"
$ws.Range("E50").Value = "  -3.13%  "

# Row 51 - RenderToken
$ws.Range("D51").Value = "'4.59"
$ws.Range("E51").Value = "  -2.44%  "
